$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '26.669.41'
$ws.Cells.Item(2, 5).Value = '  -0.08%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.597.94'
$ws.Cells.Item(3, 5).Value = '  +0.01%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.22%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''211.53'
$ws.Cells.Item(5, 5).Value = '  +0.03%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  +0.73%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.03%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +0.56%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''19.58'
$ws.Cells.Item(10, 5).Value = '  -0.77%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''0.0839'
$ws.Cells.Item(11, 5).Value = '  +0.09%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.821.82'
$ws.Cells.Item(12, 5).Value = '  +0.01%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '1.599.89'
$ws.Cells.Item(13, 5).Value = '  +0.22%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  -0.13%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '''0.524'
$ws.Cells.Item(15, 5).Value = '  +0.18%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '''65.18'
$ws.Cells.Item(16, 5).Value = '  +0.05%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '26.650.20'
$ws.Cells.Item(17, 5).Value = '  -0.09%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  +1.35%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '''210.03'
$ws.Cells.Item(19, 5).Value = '  +0.07%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +0.16%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''7.02'
$ws.Cells.Item(21, 5).Value = '  +4.06%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +0.47%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''2.33'
$ws.Cells.Item(23, 5).Value = '  +1.17%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''9.01'
$ws.Cells.Item(24, 5).Value = '  +0.95%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''144.67'
$ws.Cells.Item(25, 5).Value = '  -1.36%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.13%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''7.12'
$ws.Cells.Item(27, 5).Value = '  -0.71%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -0.64%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '''15.29'

# Row 30
$ws.Cells.Item(30, 4).Value = '''0.0516'
$ws.Cells.Item(30, 5).Value = '  +2.33%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +0.35%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''3.25'
$ws.Cells.Item(32, 5).Value = '  +0.96%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''2.96'
$ws.Cells.Item(33, 5).Value = '  +1.45%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '1.284.63'
$ws.Cells.Item(34, 5).Value = '  -0.98%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''0.621'
$ws.Cells.Item(35, 5).Value = '  -6.67%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +0.60%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +0.97%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -0.73%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +19.49%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +2.20%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -0.30%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''0.786'
$ws.Cells.Item(43, 5).Value = '  -0.54%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''63.76'
$ws.Cells.Item(44, 5).Value = '  -0.05%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '1.735.28'
$ws.Cells.Item(45, 5).Value = '  +0.07%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +0.51%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''1.57'
$ws.Cells.Item(47, 5).Value = '  -3.18%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(48, 4).Value = '0.0₆0104'
$ws.Cells.Item(48, 5).Value = '  -0.82%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Algorand'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(49, 4).Value = '''0.102'
$ws.Cells.Item(49, 5).Value = '  +1.12%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'Cronos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(50, 4).Value = '''0.0509'
$ws.Cells.Item(50, 5).Value = '  +0.72%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'USDD'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Cells.Item(51, 4).Value = '''1.00'
$ws.Cells.Item(51, 5).Value = '  -0.01%  '
